$d = $word.ActiveDocument

$wdParagraph = 4  # WdUnits.wdParagraph

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$r = $d.Content
$found = $r.Find.Execute("Docente(s)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand the found range to cover the whole heading paragraph.
    [void]$r.Expand($wdParagraph)

    # Remember which paragraph index this is, so we can fetch the sibling
    # paragraph that gets inserted right after it.
    $targetIndex = 0
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Start -eq $r.Start) {
            $targetIndex = $i
        }
    }

    # Insert a brand-new paragraph immediately after the heading.
    [void]$r.InsertParagraphAfter()

    # Turn the freshly-inserted (currently empty) paragraph into the new
    # "List Bullet" entry naming the responsible professor.
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Style = "List Bullet"
    $newPara.Range.Text = "5840917 - Fabricio Maciel Gomes"
}
